# "fixed script and schema definition in Excel"
# Normalises the camelCase column-name identifiers in the `routes` schema
# sheet to lowercase, fixes the approveRouteType column's SQL type, adds
# the missing "person" row label, and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - approveRouteType -> approveroutetype, tinyint(1) -> int(2)
$ws.Range("A3").Value = "approveroutetype"
$ws.Range("C3").Value = "int(2)"

# Row 4 - person row was missing its Japanese column-name label
$ws.Range("B4").Value = "申請者"

# Rows 6-16 - lowercase the camelCase identifiers
$ws.Range("A6").Value = "moneycondition"
$ws.Range("A7").Value = "ratecondition"
$ws.Range("A8").Value = "conditionflg"
$ws.Range("A9").Value = "approverlayer"
$ws.Range("A10").Value = "approverdept"
$ws.Range("A11").Value = "approvertitle"
$ws.Range("A12").Value = "approverid"
$ws.Range("A13").Value = "agentflg"
$ws.Range("A14").Value = "jumgflg"
$ws.Range("A16").Value = "deletereason"

# Move the active selection from A12 to B5
$ws.Range("B5").Select() | Out-Null
